# Adds a new "Socks in the Dark" heading block right after the
# "...the parrot fly across" paragraph, and relocates the trailing
# _GoBack bookmark onto the new, final (blank) paragraph of that block.
#
# Target shape (matches the OOXML diff):
#
#   ...the parrot fly across</w:r>
#   </w:p>
#   <w:p/>
#   <w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>
#   <w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
#     <w:r><w:rPr><w:b/></w:rPr><w:t>Socks in the Dark</w:t></w:r>
#   </w:p>
#   <w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>
#   <w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
#     <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#     <w:bookmarkEnd w:id="0"/>
#   </w:p>

function New-WordPackageXml($bodyFragment) {
    # Minimal FlatOPC "pkg:package" shell that Range.InsertXML expects.
    return '<?xml version="1.0"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

$d = $word.ActiveDocument

# The _GoBack bookmark currently sits at the very end of the
# "...the parrot fly across" paragraph. Drop it now; it is re-created
# on the last of the five new paragraphs below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# Find the paragraph that ends with "the parrot fly across" - it is the
# anchor the new block gets inserted after.
$hostParagraph = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*the parrot fly across*") {
        $hostParagraph = $candidate
    }
}

# The five new paragraphs, in order, as literal WordprocessingML.
$newParagraphsXml = @(
    '<w:p/>',
    '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>',
    '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Socks in the Dark</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>',
    '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
)

$previousParagraph = $hostParagraph
foreach ($paragraphXml in $newParagraphsXml) {
    # Create a fresh paragraph break after the previous one ...
    $previousParagraph.Range.InsertParagraphAfter() | Out-Null
    $newParagraph = $previousParagraph.Next()
    # ... then stamp its exact contents (pPr/runs) via InsertXML, which
    # replaces only what sits inside that paragraph's own end mark.
    $newParagraph.Range.InsertXML((New-WordPackageXml $paragraphXml)) | Out-Null
    $previousParagraph = $newParagraph
}
